# Generate Report for Handoff
#
# This script updates the localization-status workbook to reflect a new
# handoff run: the source file's GUID-based name changes from
# "dea3305b-635d-45fa-af3a-19e14bc3d44d" to
# "69e46cb5-0f80-450b-b3c7-217696f21d1c", new handoff xliff files are
# generated (new content hash "60638c12ab2af4dc1f832643049ecd8799dc72b9")
# with fresh timestamps, and since this is a fresh handoff (not yet
# handed back), the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns on the per-locale sheets are reset
# (cleared / zero-date) and their hyperlink removed.

$wb = $excel.ActiveWorkbook

$oldGuid = "dea3305b-635d-45fa-af3a-19e14bc3d44d"
$newGuid = "69e46cb5-0f80-450b-b3c7-217696f21d1c"
$newHash = "60638c12ab2af4dc1f832643049ecd8799dc72b9"

$newHoDate        = "2016-08-24 23:00:36"
$newZhHandoffDate = "2016-08-24 23:00:31"
$zeroDate         = "0001-01-01 00:00:00"

function Update-HyperlinkDisplay($ws, $cellAddress, $newDisplay) {
    # Update the hyperlink that already lives on $cellAddress in place, so
    # that its relationship id / style are preserved exactly (re-pointing
    # Address to its own current value forces the engine to refresh the
    # existing entry instead of appending a duplicate one).
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range().Address() -eq $cellAddress) {
            $addr = $hl.Address()
            $hl.Address = $addr
            $hl.TextToDisplay = $newDisplay
        }
    }
}

function Remove-Hyperlink($ws, $cellAddress) {
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range().Address() -eq $cellAddress) {
            $hl.Delete()
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = ($newGuid + ".md")
$wsOverview.Range("B2").Value = ("e2e\" + $newGuid + ".md")
$wsOverview.Range("G2").Value = $newHoDate

Update-HyperlinkDisplay $wsOverview "`$B`$2" ("e2e\" + $newGuid + ".md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = ($newGuid + ".md")
$wsZh.Range("G2").Value = ($newGuid + "." + $newHash + ".zh-cn.xlf")
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $zeroDate

Update-HyperlinkDisplay $wsZh "`$A`$2" ($newGuid + ".md")
Remove-Hyperlink $wsZh "`$I`$2"

# Reset I2's style from HyperLink back to Normal; keep J2 (which already
# had no special style) present as an explicit blank cell too.
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Style = "Normal"

$wsZh.Columns.Item(9).ColumnWidth = (18.6506053379604 - 5/6)
$wsZh.Columns.Item(10).ColumnWidth = (21.7054770333426 - 5/6)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = ($newGuid + ".md")
$wsDe.Range("G2").Value = ($newGuid + "." + $newHash + ".de-de.xlf")
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $zeroDate

Update-HyperlinkDisplay $wsDe "`$A`$2" ($newGuid + ".md")
Remove-Hyperlink $wsDe "`$I`$2"

$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Style = "Normal"

$wsDe.Columns.Item(9).ColumnWidth = (18.6506053379604 - 5/6)
$wsDe.Columns.Item(10).ColumnWidth = (21.7054770333426 - 5/6)
